$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "according to the population census data" note (row 2) is dropped;
# rows below it shift up by one.
$ws.Rows.Item(2).Delete()

# Only the 2014 figures are kept; the 1989 and 2002 columns (B and C) are
# removed, so the old column D (2014) becomes the new column B.
$ws.Range("B:C").EntireColumn.Delete()

# After the shift, row 4 holds the year header (was column D = 2014, now
# column B) and row 5 the "Area" value (803.2, carried over unchanged).
$ws.Range("B4").Value = 2014

# Match the row heights used by the regenerated export (20.1pt, custom) for
# every row, including two new trailing blank rows.
$ws.Rows.Item(1).RowHeight = 20.1
$ws.Rows.Item(2).RowHeight = 20.1
$ws.Rows.Item(3).RowHeight = 20.1
$ws.Rows.Item(4).RowHeight = 20.1
$ws.Rows.Item(5).RowHeight = 20.1
$ws.Rows.Item(6).RowHeight = 20.1
$ws.Rows.Item(7).RowHeight = 20.1
